$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cited_by_count for rows 2 and 3
$ws.Range("M2").Value = "26"
$ws.Range("M3").Value = "6"

# Rows 4 and 5 swap their article data (columns A-Q). Use Copy so that
# text-looking values (like the publication_date strings) are carried
# over verbatim instead of being re-parsed (which would turn them into
# date serial numbers).
$ws.Range("A4:Q4").Copy($ws.Range("A100:Q100"))
$ws.Range("A5:Q5").Copy($ws.Range("A4:Q4"))
$ws.Range("A100:Q100").Copy($ws.Range("A5:Q5"))
$ws.Range("A100:Q100").Clear()

# After the swap, both rows keep a cited_by_count of 4.
$ws.Range("M4").Value = "4"
$ws.Range("M5").Value = "4"
